$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 10 cell values that actually changed text ---
# "kolam" group: replace "/kolam/:id" with "/kolam/{kolam_id}" for the detail/update/destroy rows
$ws.Range("C7").Value = "/kolam/{kolam_id}"
$ws.Range("C8").Value = "/kolam/{kolam_id}"
$ws.Range("C9").Value = "/kolam/{kolam_id}"

# "pemberian pangan" group renamed to "pemberian pakan", endpoints renamed to /riwayatpakan*
$ws.Range("A11").Value = "pemberian pakan"
$ws.Range("C11").Value = "/riwayatpakan"
$ws.Range("C12").Value = "/riwayatpakan"
$ws.Range("C13").Value = "/riwayatpakan/{logpakan_id}"
$ws.Range("C14").Value = "/riwayatpakan/{logpakan_id}"
$ws.Range("C15").Value = "/riwayatpakan/{logpakan_id}"

# "loguser" group: detail endpoint renamed
$ws.Range("C18").Value = "/loguser/{loguser_id}"

# --- Column widths (closest achievable values given engine's rounding) ---
$ws.Columns.Item(1).ColumnWidth = 19.333333333333332
$ws.Columns.Item(3).ColumnWidth = 29.5
$ws.Columns.Item(5).ColumnWidth = 64

# --- Selection ---
$ws.Range("I12").Select()
